$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("UnitCombat Weaponry")

# Update footer description text
$ws.Cells.Item(11, 1).Value = "Weaponry implemented as new XML similar to promotions that are automatically gained when the unit enters a city with the appropriate bonus"

# Rename weapon/material rows, add clarifying text in parentheses
$ws.Cells.Item(3, 1).Value = "Bronze (Copper)"
$ws.Cells.Item(4, 1).Value = "Iron (Iron)"
$ws.Cells.Item(5, 1).Value = "Obsidian (Obsidian)"
$ws.Cells.Item(6, 1).Value = "Steel (Steel)"
$ws.Cells.Item(7, 1).Value = "Ash (Prime Timber)"

# Update "best weapon" value for Archer on Ash/Wood row, clear Melee/Mounted columns
$ws.Cells.Item(7, 2).Value = 2
$ws.Cells.Item(7, 3).Clear()
$ws.Cells.Item(7, 4).Clear()

# Widen column A to fit new labels
$ws.Columns.Item(1).ColumnWidth = 17.8

# Update selection to reflect where the edit focus was
$ws.Range("D8").Select()
